$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reporting period cell (B1) with the new date range
$ws.Range("B1").Value = "2021/4/14-2020/4/21"

# Move the active selection to B2 (matches the saved selection state)
$ws.Range("B2").Select()
